$wb = $excel.ActiveWorkbook
$ws = $excel.ActiveSheet

# 1) New header string for column E (write order matters for the shared
#    string table, which is built in first-seen order).
$ws.Range("E1").Value = "expected pixel size"

# 2) "camera specs" block labels (order-sensitive, see above).
$ws.Range("A9").Value = "camera specs"
$ws.Range("B10").Value = "image width [px]"
$ws.Range("B11").Value = "FOV [rad]"

# 3) Column E formulas (expected pixel size). These reference A10/A11
#    which are still blank at this point; that's fine, they get fixed up
#    once A10/A11 receive real values below. Written before touching
#    column D so this new shared-formula group is allocated before D's
#    existing group is rewritten.
$ws.Range("E2").Formula = "=A2*`$A`$10/(2*B2*TAN(`$A`$11/2))"
$ws.Range("E3:E5").Formula = "=A3*`$A`$10/(2*B3*TAN(`$A`$11/2))"

# 4) Column D formulas: swap the old C2*180/PI() idiom for DEGREES().
$ws.Range("D2").Formula = "=DEGREES(C2)"
$ws.Range("D3:D5").Formula = "=DEGREES(C3)"

# 5) Fill in the actual camera-spec values (this fixes up column E's
#    results now that A10/A11 hold real numbers).
$ws.Range("A10").Value = 640
$ws.Range("A11").Formula = "=RADIANS(87)"

# 6) New column F: actual pixel size counted from the image (plain
#    numbers, written last so its header claims the final shared-string
#    slot).
$ws.Range("F1").Value = "actual pixel size (counted from image)"
$ws.Range("F2").Value = 11.5
$ws.Range("F3").Value = 11
$ws.Range("F4").Value = 11.5
$ws.Range("F5").Value = 12.5

# --- Reposition the two charts to make room for the new columns ---
$co1 = $ws.ChartObjects(1)
$co1.Left = 419.5625
$co1.Top = 16.12496062992126
$co1.Width = 433.0625
$co1.Height = 216.0

$co2 = $ws.ChartObjects(2)
$co2.Left = 880.6874606299212
$co2.Top = 17.62496062992126
$co2.Width = 433.0625
$co2.Height = 216.0

# --- Selection moves to F9 ---
[void]$ws.Range("F9").Select()
